# Update column F ("dSF") values on the active sheet to reflect the
# repulled/recalculated data from the commit "repull data, push all data, mean calculation".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = 7
    8  = 3
    9  = -4
    10 = -6
    13 = -9
    15 = -2
    18 = -4
    19 = -4
    21 = -5
    22 = -2
    24 = -5
    37 = -6
    38 = 1
    43 = -2
    48 = 3
    50 = 0
    53 = 0
    55 = -6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
